$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.928.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.498.53'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.85'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.28'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.497.92'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.44%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.097.51'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000179'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.86%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.499.79'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.033.86'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.00'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.53'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.66'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.27'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.639.10'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.07'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.43%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -7.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.25'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.24'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.505.27'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.87'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.145'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.94%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '167.72'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0809'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.08'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.812'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.97'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.09%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.39'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.89'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.438.04'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.898'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.69%  '
